$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("summary")

$ws1.Cells.Item(2, 4).Value = 1985
$ws1.Cells.Item(2, 5).Value = 81.21
$ws1.Cells.Item(2, 6).Value = -1.78
$ws1.Cells.Item(2, 8).Value = 0.99
$ws1.Cells.Item(2, 9).Value = -0.13
$ws1.Cells.Item(2, 10).Value = 0.21
$ws1.Cells.Item(2, 11).Value = 0.06
$ws1.Cells.Item(2, 12).Value = 1.16
$ws1.Cells.Item(3, 4).Value = 1862
$ws1.Cells.Item(3, 5).Value = 49.73
$ws1.Cells.Item(3, 6).Value = 0.03
$ws1.Cells.Item(3, 9).Value = 0.04
$ws1.Cells.Item(3, 10).Value = 0.22
$ws1.Cells.Item(4, 4).Value = 1379
$ws1.Cells.Item(4, 5).Value = 18.56
$ws1.Cells.Item(4, 6).Value = 1.81
$ws1.Cells.Item(4, 7).Value = 0.07
$ws1.Cells.Item(4, 8).Value = 1.02
$ws1.Cells.Item(4, 9).Value = 0.59
$ws1.Cells.Item(4, 12).Value = 0.89
$ws1.Cells.Item(5, 3).Value = 680
$ws1.Cells.Item(5, 4).Value = 657
$ws1.Cells.Item(5, 5).Value = 75.95
$ws1.Cells.Item(5, 6).Value = -1.5
$ws1.Cells.Item(5, 9).Value = 0.42
$ws1.Cells.Item(5, 10).Value = 0.22
$ws1.Cells.Item(5, 11).Value = 0.03
$ws1.Cells.Item(5, 12).Value = 0.88
$ws1.Cells.Item(6, 3).Value = 680
$ws1.Cells.Item(6, 4).Value = 655
$ws1.Cells.Item(6, 5).Value = 72.52
$ws1.Cells.Item(6, 9).Value = 0.48
$ws1.Cells.Item(7, 3).Value = 680
$ws1.Cells.Item(7, 4).Value = 647
$ws1.Cells.Item(7, 5).Value = 64.91
$ws1.Cells.Item(7, 6).Value = -0.87
$ws1.Cells.Item(7, 8).Value = 1
$ws1.Cells.Item(7, 9).Value = 0.05
$ws1.Cells.Item(7, 10).Value = 0.28
$ws1.Cells.Item(7, 12).Value = 1.1
$ws1.Cells.Item(8, 3).Value = 680
$ws1.Cells.Item(8, 4).Value = 648
$ws1.Cells.Item(8, 5).Value = 62.35
$ws1.Cells.Item(8, 6).Value = -0.74
$ws1.Cells.Item(8, 9).Value = 0.77
$ws1.Cells.Item(8, 12).Value = 0.84
$ws1.Cells.Item(9, 3).Value = 714
$ws1.Cells.Item(9, 4).Value = 676
$ws1.Cells.Item(9, 5).Value = 63.76
$ws1.Cells.Item(9, 6).Value = -0.7
$ws1.Cells.Item(9, 8).Value = 0.98
$ws1.Cells.Item(9, 9).Value = -0.62
$ws1.Cells.Item(9, 10).Value = 0.32
$ws1.Cells.Item(9, 11).Value = 0.05
$ws1.Cells.Item(9, 12).Value = 1.3
$ws1.Cells.Item(10, 3).Value = 714
$ws1.Cells.Item(10, 4).Value = 643
$ws1.Cells.Item(10, 5).Value = 55.37
$ws1.Cells.Item(10, 6).Value = -0.27
$ws1.Cells.Item(10, 8).Value = 1.03
$ws1.Cells.Item(10, 9).Value = 0.93
$ws1.Cells.Item(10, 11).Value = 0.05
$ws1.Cells.Item(10, 12).Value = 0.84
$ws1.Cells.Item(11, 3).Value = 714
$ws1.Cells.Item(11, 4).Value = 623
$ws1.Cells.Item(11, 5).Value = 51.04
$ws1.Cells.Item(11, 6).Value = -0.07
$ws1.Cells.Item(11, 9).Value = -0.52
$ws1.Cells.Item(11, 10).Value = 0.35
$ws1.Cells.Item(11, 11).Value = 0.07
$ws1.Cells.Item(11, 12).Value = 1.2
$ws1.Cells.Item(12, 3).Value = 714
$ws1.Cells.Item(12, 5).Value = 47.62
$ws1.Cells.Item(12, 6).Value = 0.11
$ws1.Cells.Item(12, 8).Value = 1
$ws1.Cells.Item(12, 9).Value = 0.11
$ws1.Cells.Item(12, 10).Value = 0.33
$ws1.Cells.Item(12, 11).Value = 0.05
$ws1.Cells.Item(12, 12).Value = 1.02
$ws1.Cells.Item(13, 3).Value = 714
$ws1.Cells.Item(13, 5).Value = 42.47
$ws1.Cells.Item(13, 6).Value = 0.35
$ws1.Cells.Item(13, 8).Value = 1
$ws1.Cells.Item(13, 9).Value = 0.09
$ws1.Cells.Item(13, 10).Value = 0.33
$ws1.Cells.Item(13, 11).Value = 0.03
$ws1.Cells.Item(13, 12).Value = 1.06
$ws1.Cells.Item(14, 3).Value = 706
$ws1.Cells.Item(14, 4).Value = 671
$ws1.Cells.Item(14, 5).Value = 39.64
$ws1.Cells.Item(14, 6).Value = 0.64
$ws1.Cells.Item(14, 8).Value = 1
$ws1.Cells.Item(14, 9).Value = -0.04
$ws1.Cells.Item(14, 10).Value = 0.32
$ws1.Cells.Item(14, 12).Value = 1.11
$ws1.Cells.Item(15, 3).Value = 706
$ws1.Cells.Item(15, 4).Value = 678
$ws1.Cells.Item(15, 5).Value = 37.46
$ws1.Cells.Item(15, 6).Value = 0.76
$ws1.Cells.Item(15, 9).Value = 0.08
$ws1.Cells.Item(15, 12).Value = 1.09
$ws1.Cells.Item(16, 3).Value = 706
$ws1.Cells.Item(16, 4).Value = 676
$ws1.Cells.Item(16, 5).Value = 26.48
$ws1.Cells.Item(16, 6).Value = 1.38
$ws1.Cells.Item(16, 9).Value = -1.31
$ws1.Cells.Item(16, 10).Value = 0.38
$ws1.Cells.Item(16, 12).Value = 1.82
$ws1.Cells.Item(17, 3).Value = 706
$ws1.Cells.Item(17, 4).Value = 676
$ws1.Cells.Item(17, 5).Value = 23.08
$ws1.Cells.Item(17, 9).Value = -0.7
$ws1.Cells.Item(17, 10).Value = 0.33
$ws1.Cells.Item(17, 11).Value = 0.05
$ws1.Cells.Item(17, 12).Value = 1.47

$ws2 = $wb.Worksheets.Item("model_fit")

$ws2.Cells.Item(2, 2).Value = 2100
$ws2.Cells.Item(2, 4).Value = 15677
$ws2.Cells.Item(2, 5).Value = 15711
$ws2.Cells.Item(2, 6).Value = 15807
$ws2.Cells.Item(2, 7).Value = 0.552
$ws2.Cells.Item(2, 8).Value = 0.364
$ws2.Cells.Item(3, 2).Value = 2100
$ws2.Cells.Item(3, 4).Value = 15652
$ws2.Cells.Item(3, 5).Value = 15716
$ws2.Cells.Item(3, 6).Value = 15896
$ws2.Cells.Item(3, 8).Value = 0.355
